$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.327.71"
$ws.Range("E2").Value = "  +2.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.59"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -1.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.30"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("E6").Value = "  -1.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5136"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08463"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.60"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.118"
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.268"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.902.31"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.67"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.342"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.39"
$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06743"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.94"
$ws.Range("E20").Value = "  +0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.035"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.341.32"
$ws.Range("E23").Value = "  +2.86%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").Value = "  -2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.118.88"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.67"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.02"
$ws.Range("E28").Value = "  +0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.458"
$ws.Range("E29").Value = "  +2.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.22"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.064"
$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.078"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.663"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02480"
$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06594"
$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.109"
$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2204"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.239"
$ws.Range("E39").Value = "  +3.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.136"
$ws.Range("E40").Value = "  +1.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6515"
$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.232"
$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6067"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.683"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.058"
$ws.Range("E47").Value = "  +2.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.232"
$ws.Range("E48").Value = "  +1.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.38"
$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("E51").Value = "  +0.92%  "
